$wb = $excel.ActiveWorkbook

# New "want to go" counts for column F, keyed by row number.
$updates = @{
    2  = 1115
    3  = 441
    4  = 1529
    5  = 8831
    6  = 97
    7  = 497
    8  = 658
    9  = 311
    11 = 26
    12 = 32
    13 = 3700
    14 = 53
    16 = 92
    17 = 3352
    18 = 153
    19 = 1128
    20 = 321
    21 = 223
    22 = 2487
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

# Sheet "展览" also has row 23 updated (its row 24 does not exist there).
$wsA = $wb.Worksheets.Item("展览")
$wsA.Range("F23").Value = 87

# Sheet "全部类型" has the equivalent update on row 24 instead of row 23
# (row 23 there holds a different record and stays unchanged).
$wsB = $wb.Worksheets.Item("全部类型")
$wsB.Range("F24").Value = 87
